# Generate Report for Handoff
# Reorders the per-file rows (by filename) on the Overview / zh-cn / de-de
# sheets and flips the now-last row (d972d633-...) from "Handed back" to
# "Ready for handoff", refreshing its handoff/target/handback file+time
# columns to the newer values produced by the fresh handoff.

$wb = $excel.ActiveWorkbook

# ---------------- Overview sheet ----------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "ffffe2659819-143f-4037-b1c0-5a738e2da752.md"
$ws.Range("A3").Value = "ffffff7d2efb25-3a66-455f-96d4-29450528165a.md"
$ws.Range("A4").Value = "d972d633-8fdb-4f90-a626-2f1d93d543da.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

# ---------------- zh-cn sheet ----------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "ffffe2659819-143f-4037-b1c0-5a738e2da752.md"
$ws.Range("C2").Value = "5ea0682b-e9cb-451d-b591-dcb546182c0e.5d27baed4f0926059f51a0c7bb0078099d7a0bfa.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-10 14:53:44"
$ws.Range("E2").Value = "5ea0682b-e9cb-451d-b591-dcb546182c0e.md"
$ws.Range("F2").Value = "5ea0682b-e9cb-451d-b591-dcb546182c0e.5d27baed4f0926059f51a0c7bb0078099d7a0bfa.zh-cn.xlf"
$ws.Range("G2").Value = "2016-03-10 14:54:07"

$ws.Range("A3").Value = "ffffff7d2efb25-3a66-455f-96d4-29450528165a.md"

$ws.Range("A4").Value = "d972d633-8fdb-4f90-a626-2f1d93d543da.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "d972d633-8fdb-4f90-a626-2f1d93d543da.2307ed0e4ac8f83a1af97480276c46b83f03b0cd.zh-cn.xlf"
$ws.Range("D4").Value = "2016-03-10 14:55:48"
$ws.Range("E4").Value = "d972d633-8fdb-4f90-a626-2f1d93d543da.md"
$ws.Range("F4").Value = "d972d633-8fdb-4f90-a626-2f1d93d543da.2307ed0e4ac8f83a1af97480276c46b83f03b0cd.zh-cn.xlf"
$ws.Range("G4").Value = "2016-03-10 14:55:19"

# ---------------- de-de sheet ----------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "ffffe2659819-143f-4037-b1c0-5a738e2da752.md"
$ws.Range("C2").Value = "5ea0682b-e9cb-451d-b591-dcb546182c0e.5d27baed4f0926059f51a0c7bb0078099d7a0bfa.de-de.xlf"
$ws.Range("D2").Value = "2016-03-10 14:53:48"
$ws.Range("E2").Value = "5ea0682b-e9cb-451d-b591-dcb546182c0e.md"
$ws.Range("F2").Value = "5ea0682b-e9cb-451d-b591-dcb546182c0e.5d27baed4f0926059f51a0c7bb0078099d7a0bfa.de-de.xlf"
$ws.Range("G2").Value = "2016-03-10 14:54:15"

$ws.Range("A3").Value = "ffffff7d2efb25-3a66-455f-96d4-29450528165a.md"

$ws.Range("A4").Value = "d972d633-8fdb-4f90-a626-2f1d93d543da.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "d972d633-8fdb-4f90-a626-2f1d93d543da.2307ed0e4ac8f83a1af97480276c46b83f03b0cd.de-de.xlf"
$ws.Range("D4").Value = "2016-03-10 14:55:51"
$ws.Range("E4").Value = "d972d633-8fdb-4f90-a626-2f1d93d543da.md"
$ws.Range("F4").Value = "d972d633-8fdb-4f90-a626-2f1d93d543da.2307ed0e4ac8f83a1af97480276c46b83f03b0cd.de-de.xlf"
$ws.Range("G4").Value = "2016-03-10 14:55:27"
